# translations changes before pull from upstream
#
# Adds new "greek" locale row to the settings sheet and a batch of new
# beneficiary/delivery/override/search related translation keys to the
# common_translations sheet, then updates the saved selection/window state.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("common_translations")
$ws4 = $wb.Worksheets.Item("settings")
$ws5 = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------------
# settings!row10 - new "greek" language configuration row (mirrors the
# existing english/hindi/spanish rows directly above it).
# ---------------------------------------------------------------------------
$ws4.Range("A10").Value = "greek"
$ws4.Range("F10").Value = "Greek"
$ws4.Range("G10").Value = "Greek"
$ws4.Range("H10").Value = "Greek (as greek name)"

# Match the existing header-row styling (s="2") used by the other language
# rows' key cells - copy formats from the row above (A9/F9/G9/H9).
$ws4.Range("A9").Copy()
$ws4.Range("A10").PasteSpecial(-4122)
$ws4.Range("F9").Copy()
$ws4.Range("F10:G10").PasteSpecial(-4122)
$ws4.Range("H9").Copy()
$ws4.Range("H10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# common_translations - new rows 21-53 of beneficiary/delivery/search/
# override/authorization/item-pack translation keys.
# ---------------------------------------------------------------------------
$ws3.Range("A21").Value = "no_active_beneficiary"
$ws3.Range("B21").Value = "No Active Beneficiary Detected"

$ws3.Range("B22").Value = "No Disabled Beneficiary Detected"
$ws3.Range("A22").Value = "no_disabled_beneficiary"

$ws3.Range("A23").Value = "enter_beneficiary_code"
$ws3.Range("B23").Value = "Please Enter Beneficiary Code"

$ws3.Range("A24").Value = "active_beneficiaries_title"
$ws3.Range("B24").Value = "Active Beneficiaries"

$ws3.Range("B25").Value = "Disabled Beneficiaries"

$ws3.Range("A26").Value = "beneficiary_data_title"
$ws3.Range("B26").Value = "Beneficiary Data"

$ws3.Range("A27").Value = "delivery_data_title"
$ws3.Range("B27").Value = "Delivery Data"

$ws3.Range("A28").Value = "beneficiary_lists"
$ws3.Range("B28").Value = "Beneficiary Lists"

$ws3.Range("A29").Value = "view_all_deliveries"
$ws3.Range("B29").Value = "View All Deliveries"

$ws3.Range("A30").Value = "advanced_search"
$ws3.Range("B30").Value = "Advanced Search"

$ws3.Range("A31").Value = "data_categories"
$ws3.Range("A32").Value = "view_beneficiary_data"
$ws3.Range("A33").Value = "view_delivery_data"
$ws3.Range("B31").Value = "Data Categories"
$ws3.Range("B33").Value = "View Delivery Data"
$ws3.Range("B32").Value = "View Beneficiary Data"

$ws3.Range("A34").Value = "enable_beneficiary"
$ws3.Range("B34").Value = "Enable Beneficiary"

$ws3.Range("A35").Value = "disable_beneficiary"
$ws3.Range("B35").Value = "Disable Beneficiary"

$ws3.Range("A36").Value = "choose_method"
$ws3.Range("B36").Value = "Choose Method"

$ws3.Range("B37").Value = "Override Options"
$ws3.Range("A37").Value = "override_options"

$ws3.Range("A38").Value = "override_registration"
$ws3.Range("B38").Value = "Override Registration"

$ws3.Range("A39").Value = "override_entitlement"
$ws3.Range("B39").Value = "Override Entitlement"

$ws3.Range("A40").Value = "serach_beneficiaries"
$ws3.Range("B40").Value = "Search for Beneficiaries"

$ws3.Range("A41").Value = "search_deliveries"
$ws3.Range("B41").Value = "Search for Deliveries"

$ws3.Range("A42").Value = "beneficiary"
$ws3.Range("B42").Value = "beneficiary"

$ws3.Range("A43").Value = "beneficiaries"
$ws3.Range("B43").Value = "beneficiaries"

$ws3.Range("A44").Value = "delivery"
$ws3.Range("B44").Value = "delivery"

$ws3.Range("A45").Value = "deliveries"
$ws3.Range("B45").Value = "deliveries"

$ws3.Range("A46").Value = "found"
$ws3.Range("B46").Value = "found"

$ws3.Range("A47").Value = "invalid_search"
$ws3.Range("B47").Value = "Invalid Search"

$ws3.Range("A48").Value = "authorization_name"
$ws3.Range("B48").Value = "Authorization Name"

$ws3.Range("A49").Value = "authorization_id"
$ws3.Range("B49").Value = "Authorization ID"

$ws3.Range("A50").Value = "item_pack_name"
$ws3.Range("B50").Value = "Item Pack Name"

$ws3.Range("A51").Value = "item_pack_description"
$ws3.Range("B51").Value = "Item Pack Description"

$ws3.Range("A52").Value = "item_pack_id"
$ws3.Range("B52").Value = "Item Pack ID"

$ws3.Range("A53").Value = "beneficiary_code"
$ws3.Range("B53").Value = "Beneficiary Code"

# Written last so this key lands on the final shared-string slot, matching
# the upstream edit (the author filled this cell in after everything else).
$ws3.Range("A25").Value = "disabled_beneficiaries_title"

# Row 22 ("no disabled beneficiary detected") wraps onto two lines like the
# other long prompts in this sheet (rows 11/15), so it needs the taller
# row height.
$ws3.Rows.Item(22).RowHeight = 26

# Rows 48-53 (authorization / item pack / beneficiary code) pick up the
# alternate "choices"-style formatting (style index 3 in the original
# workbook) that's already used elsewhere in the workbook - copy it over.
$ws5.Range("A2").Copy()
$ws3.Range("A48:C52").PasteSpecial(-4122)
$ws5.Range("A2").Copy()
$ws3.Range("A53:B53").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Selection / window bookkeeping, mirroring the saved-state portion of the
# diff.
# ---------------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("A25").Select()

$ws4.Activate()
$ws4.Range("G14").Select()

$ws3.Activate()

$aw = $excel.ActiveWindow
$aw.Left = 0
$aw.Top = 460
$aw.Width = 12020
$aw.Height = 16220
